# preparation publication 0.2.0
# - bump Version to 0.2.0
# - bump Date to 2023-10-20T08:59:58+00:00
# - add a new "Jurisdiction" / "iso:code:3166:FR" row right after "Contact"
#   (this pushes Description/Purpose/Copyright/... etc. down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 is currently "Description"; insert a new blank row there so that
# "Description" (and everything below it) moves down to row 12.
$ws.Rows.Item(11).Insert()

# The freshly inserted row doesn't carry the table's normal bordered /
# wrap-text style, so copy formatting from the row right below (which still
# has the correct look) onto the new row.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new Jurisdiction row.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Bump the Version and Date metadata values.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"
